$d = $word.ActiveDocument

# The site footer that Jekyll stamps onto every rebuilt page was dropped:
#   "Ver no Jupiter Salvar em pdf Salvar em docx"
#   "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages.
#    Original theme under Creative Commons Attribution"
# plus the blank paragraph that used to sit between that footer and the
# trailing page-break paragraph. Locate the block by text (robust to any
# shift in paragraph indices) and delete it as a single contiguous range,
# leaving the blank paragraph right after "LOB1206: Solos I (Requisito
# fraco)" and the page-break paragraph untouched.

$wdParagraph = 4

$start = $d.Content.Duplicate
$foundStart = $start.Find.Execute(
    "Ver no Jupiter Salvar em pdf Salvar em docx",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
[void]$start.Expand($wdParagraph)

$end = $d.Content.Duplicate
$foundEnd = $end.Find.Execute(
    "Powered by Jekyll and Github pages",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
[void]$end.Expand($wdParagraph)

# Pull in the blank separator paragraph immediately following the copyright
# line so it is removed along with the rest of the footer block.
$sep = $d.Range($end.End, $end.End)
[void]$sep.Expand($wdParagraph)

$target = $d.Range($start.Start, $sep.End)
[void]$target.Delete()
